$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20/21, 31/32, 44/45 full-row swaps (coin name, link, price, volume) ---
# Row 20/21 swap: ShibaInu <-> InternetComputer(DFINITY)
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").Value = "'13.57"
$ws.Range("E20").Value = "  +6.64%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "'0.0000104"
$ws.Range("E21").Value = "  -0.28%  "

# Row 31/32 swap: InjectiveProtocol <-> Monero
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'165.10"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("B32").Value = "InjectiveProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D32").Value = "'35.58"
$ws.Range("E32").Value = "  -4.13%  "

# Row 44/45 swap: Algorand <-> FirstDigitalUSD
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.01"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.226"
$ws.Range("E45").Value = "  +1.10%  "

# --- Remaining single-cell updates (price / volume columns) ---
$ws.Range("D2").Value = "42.555.58"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "2.292.91"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").Value = "'311.66"
$ws.Range("E5").Value = "  -2.12%  "
$ws.Range("D6").Value = "'104.03"
$ws.Range("E6").Value = "  +3.16%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'39.11"
$ws.Range("E10").Value = "  +0.63%  "
$ws.Range("E11").Value = "  +0.43%  "
$ws.Range("D12").Value = "'8.27"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "'0.984"
$ws.Range("E14").Value = "  +3.59%  "
$ws.Range("D15").Value = "'15.10"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("D16").Value = "2.642.61"
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.299.84"
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "42.727.19"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("D22").Value = "'73.14"
$ws.Range("E22").Value = "  +0.87%  "
$ws.Range("D23").Value = "'3.44"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("D24").Value = "'263.06"
$ws.Range("E24").Value = "  -1.19%  "
$ws.Range("D25").Value = "'2.18"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "'10.70"
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("D28").Value = "'7.04"
$ws.Range("E28").Value = "  +16.58%  "
$ws.Range("D29").Value = "'2.26"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").Value = "'22.28"
$ws.Range("E30").Value = "  -0.49%  "
$ws.Range("D33").Value = "'0.0860"
$ws.Range("E33").Value = "  -0.83%  "
$ws.Range("D34").Value = "'0.130"
$ws.Range("E34").Value = "  -1.78%  "
$ws.Range("E35").Value = "  +0.16%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("E36").Value = "  -2.11%  "
$ws.Range("D37").Value = "'4.50"
$ws.Range("E37").Value = "  -1.22%  "
$ws.Range("D38").Value = "'0.0349"
$ws.Range("E38").Value = "  -1.13%  "
$ws.Range("D39").Value = "'3.72"
$ws.Range("E39").Value = "  +2.27%  "
$ws.Range("D40").Value = "'2.72"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'1.58"
$ws.Range("E41").Value = "  +5.12%  "
$ws.Range("D42").Value = "'99.06"
$ws.Range("E42").Value = "  +8.68%  "
$ws.Range("D43").Value = "'69.14"
$ws.Range("E43").Value = "  +1.61%  "
$ws.Range("D46").Value = "1.729.56"
$ws.Range("E46").Value = "  +7.68%  "
$ws.Range("D47").Value = "'11.97"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("D48").Value = "'78.63"
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'110.64"
$ws.Range("E49").Value = "  -3.43%  "
$ws.Range("E50").Value = "  -0.19%  "
$ws.Range("D51").Value = "'8.65"
$ws.Range("E51").Value = "  -3.02%  "
